$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---- Update cell values (shift Y:AE columns left by one; adjust column A target list; insert storeKeys into json list) ----
$ws.Cells.Item(1, 27).Value = 'webcookie'
$ws.Cells.Item(1, 28).Value = 'ws'
$ws.Cells.Item(1, 29).Value = 'ws.async'
$ws.Cells.Item(1, 30).Value = 'xml'
$ws.Cells.Item(1, 25).Value = 'web'
$ws.Cells.Item(1, 26).Value = 'webalert'
$ws.Cells.Item(2, 27).Value = 'assertNotPresent(name)'
$ws.Cells.Item(2, 28).Value = 'assertReturnCode(var,returnCode)'
$ws.Cells.Item(2, 29).Value = 'delete(url,body,output)'
$ws.Cells.Item(2, 30).Value = 'append(xml,xpath,content,var)'
$ws.Cells.Item(2, 25).Value = 'assertAndClick(locator,label)'
$ws.Cells.Item(2, 26).Value = 'accept()'
$ws.Cells.Item(3, 27).Value = 'assertPresent(name)'
$ws.Cells.Item(3, 28).Value = 'delete(url,body,var)'
$ws.Cells.Item(3, 29).Value = 'download(url,queryString,saveTo)'
$ws.Cells.Item(3, 30).Value = 'assertCorrectness(xml,schema)'
$ws.Cells.Item(3, 25).Value = 'assertAttribute(locator,attrName,value)'
$ws.Cells.Item(3, 26).Value = 'assertPresent()'
$ws.Cells.Item(4, 27).Value = 'assertValue(name,value)'
$ws.Cells.Item(4, 28).Value = 'download(url,queryString,saveTo)'
$ws.Cells.Item(4, 29).Value = 'get(url,queryString,output)'
$ws.Cells.Item(4, 30).Value = 'assertElementCount(xml,xpath,count)'
$ws.Cells.Item(4, 25).Value = 'assertAttributeContains(locator,attrName,contains)'
$ws.Cells.Item(4, 26).Value = 'assertText(text,matchBy)'
$ws.Cells.Item(5, 27).Value = 'delete(name)'
$ws.Cells.Item(5, 28).Value = 'get(url,queryString,var)'
$ws.Cells.Item(5, 29).Value = 'head(url,output)'
$ws.Cells.Item(5, 30).Value = 'assertElementNotPresent(xml,xpath)'
$ws.Cells.Item(5, 25).Value = 'assertAttributeNotContains(locator,attrName,contains)'
$ws.Cells.Item(5, 26).Value = 'dismiss()'
$ws.Cells.Item(6, 27).Value = 'deleteAll()'
$ws.Cells.Item(6, 28).Value = 'head(url,var)'
$ws.Cells.Item(6, 29).Value = 'patch(url,body,output)'
$ws.Cells.Item(6, 30).Value = 'assertElementPresent(xml,xpath)'
$ws.Cells.Item(6, 25).Value = 'assertAttributeNotPresent(locator,attrName)'
$ws.Cells.Item(6, 26).Value = 'replyCancel(text)'
$ws.Cells.Item(7, 27).Value = 'save(var,name)'
$ws.Cells.Item(7, 28).Value = 'header(name,value)'
$ws.Cells.Item(7, 29).Value = 'post(url,body,output)'
$ws.Cells.Item(7, 30).Value = 'assertSoap(wsdl,xml)'
$ws.Cells.Item(7, 25).Value = 'assertAttributePresent(locator,attrName)'
$ws.Cells.Item(7, 26).Value = 'replyOK(text)'
$ws.Cells.Item(8, 27).Value = 'saveAll(var)'
$ws.Cells.Item(8, 28).Value = 'headerByVar(name,var)'
$ws.Cells.Item(8, 29).Value = 'put(url,body,output)'
$ws.Cells.Item(8, 30).Value = 'assertSoapFaultCode(expected,xml)'
$ws.Cells.Item(8, 25).Value = 'assertChecked(locator)'
$ws.Cells.Item(8, 26).Value = 'storeText(var)'
$ws.Cells.Item(9, 28).Value = 'jwtParse(var,token,key)'
$ws.Cells.Item(9, 30).Value = 'assertSoapFaultString(expected,xml)'
$ws.Cells.Item(9, 25).Value = 'assertContainCount(locator,text,count)'
$ws.Cells.Item(10, 28).Value = 'jwtSignHS256(var,payload,key)'
$ws.Cells.Item(10, 30).Value = 'assertValue(xml,xpath,expected)'
$ws.Cells.Item(10, 25).Value = 'assertCssNotPresent(locator,property)'
$ws.Cells.Item(11, 28).Value = 'oauth(var,url,auth)'
$ws.Cells.Item(11, 30).Value = 'assertValues(xml,xpath,array,exactOrder)'
$ws.Cells.Item(11, 25).Value = 'assertCssPresent(locator,property,value)'
$ws.Cells.Item(12, 28).Value = 'patch(url,body,var)'
$ws.Cells.Item(12, 30).Value = 'assertWellformed(xml)'
$ws.Cells.Item(12, 25).Value = 'assertElementByAttributes(nameValues)'
$ws.Cells.Item(13, 28).Value = 'post(url,body,var)'
$ws.Cells.Item(13, 30).Value = 'beautify(xml,var)'
$ws.Cells.Item(13, 25).Value = 'assertElementByText(locator,text)'
$ws.Cells.Item(14, 28).Value = 'put(url,body,var)'
$ws.Cells.Item(14, 30).Value = 'clear(xml,xpath,var)'
$ws.Cells.Item(14, 25).Value = 'assertElementCount(locator,count)'
$ws.Cells.Item(15, 28).Value = 'saveResponsePayload(var,file,append)'
$ws.Cells.Item(15, 30).Value = 'delete(xml,xpath,var)'
$ws.Cells.Item(15, 25).Value = 'assertElementNotPresent(locator)'
$ws.Cells.Item(16, 28).Value = 'soap(action,url,payload,var)'
$ws.Cells.Item(16, 30).Value = 'insertAfter(xml,xpath,content,var)'
$ws.Cells.Item(16, 13).Value = 'storeKeys(json,jsonpath,var)'
$ws.Cells.Item(16, 25).Value = 'assertElementPresent(locator)'
$ws.Cells.Item(17, 28).Value = 'upload(url,body,fileParams,var)'
$ws.Cells.Item(17, 30).Value = 'insertBefore(xml,xpath,content,var)'
$ws.Cells.Item(17, 13).Value = 'storeValue(json,jsonpath,var)'
$ws.Cells.Item(17, 25).Value = 'assertElementsPresent(prefix)'
$ws.Cells.Item(18, 30).Value = 'minify(xml,var)'
$ws.Cells.Item(18, 13).Value = 'storeValues(json,jsonpath,var)'
$ws.Cells.Item(18, 25).Value = 'assertFocus(locator)'
$ws.Cells.Item(19, 30).Value = 'prepend(xml,xpath,content,var)'
$ws.Cells.Item(19, 25).Value = 'assertFrameCount(count)'
$ws.Cells.Item(20, 30).Value = 'replace(xml,xpath,content,var)'
$ws.Cells.Item(20, 25).Value = 'assertFramePresent(frameName)'
$ws.Cells.Item(21, 30).Value = 'replaceIn(xml,xpath,content,var)'
$ws.Cells.Item(21, 25).Value = 'assertIECompatMode()'
$ws.Cells.Item(22, 30).Value = 'storeCount(xml,xpath,var)'
$ws.Cells.Item(22, 25).Value = 'assertIENativeMode()'
$ws.Cells.Item(23, 30).Value = 'storeSoapFaultCode(var,xml)'
$ws.Cells.Item(23, 25).Value = 'assertLinkByLabel(label)'
$ws.Cells.Item(24, 30).Value = 'storeSoapFaultDetail(var,xml)'
$ws.Cells.Item(24, 25).Value = 'assertNotChecked(locator)'
$ws.Cells.Item(25, 1).Value = 'web'
$ws.Cells.Item(25, 30).Value = 'storeSoapFaultString(var,xml)'
$ws.Cells.Item(25, 25).Value = 'assertNotFocus(locator)'
$ws.Cells.Item(26, 1).Value = 'webalert'
$ws.Cells.Item(26, 30).Value = 'storeValue(xml,xpath,var)'
$ws.Cells.Item(26, 25).Value = 'assertNotText(locator,text)'
$ws.Cells.Item(27, 1).Value = 'webcookie'
$ws.Cells.Item(27, 30).Value = 'storeValues(xml,xpath,var)'
$ws.Cells.Item(27, 25).Value = 'assertNotVisible(locator)'
$ws.Cells.Item(28, 1).Value = 'ws'
$ws.Cells.Item(28, 25).Value = 'assertOneMatch(locator)'
$ws.Cells.Item(29, 1).Value = 'ws.async'
$ws.Cells.Item(29, 25).Value = 'assertScrollbarHNotPresent(locator)'
$ws.Cells.Item(30, 1).Value = 'xml'
$ws.Cells.Item(30, 25).Value = 'assertScrollbarHPresent(locator)'
$ws.Cells.Item(31, 25).Value = 'assertScrollbarVNotPresent(locator)'
$ws.Cells.Item(32, 25).Value = 'assertScrollbarVPresent(locator)'
$ws.Cells.Item(33, 25).Value = 'assertTable(locator,row,column,text)'
$ws.Cells.Item(34, 25).Value = 'assertText(locator,text)'
$ws.Cells.Item(35, 25).Value = 'assertTextContains(locator,text)'
$ws.Cells.Item(36, 25).Value = 'assertTextCount(locator,text,count)'
$ws.Cells.Item(37, 25).Value = 'assertTextList(locator,list,ignoreOrder)'
$ws.Cells.Item(38, 25).Value = 'assertTextMatches(text,minMatch,scrollTo)'
$ws.Cells.Item(39, 25).Value = 'assertTextNotContains(locator,text)'
$ws.Cells.Item(40, 25).Value = 'assertTextNotPresent(text)'
$ws.Cells.Item(41, 25).Value = 'assertTextOrder(locator,descending)'
$ws.Cells.Item(42, 25).Value = 'assertTextPresent(text)'
$ws.Cells.Item(43, 25).Value = 'assertTitle(text)'
$ws.Cells.Item(44, 25).Value = 'assertValue(locator,value)'
$ws.Cells.Item(45, 25).Value = 'assertValueOrder(locator,descending)'
$ws.Cells.Item(46, 25).Value = 'assertVisible(locator)'
$ws.Cells.Item(47, 25).Value = 'checkAll(locator)'
$ws.Cells.Item(48, 25).Value = 'clearLocalStorage()'
$ws.Cells.Item(49, 25).Value = 'click(locator)'
$ws.Cells.Item(50, 25).Value = 'clickAll(locator)'
$ws.Cells.Item(51, 25).Value = 'clickAndWait(locator,waitMs)'
$ws.Cells.Item(52, 25).Value = 'clickByLabel(label)'
$ws.Cells.Item(53, 25).Value = 'clickByLabelAndWait(label,waitMs)'
$ws.Cells.Item(54, 25).Value = 'clickOffset(locator,x,y)'
$ws.Cells.Item(55, 25).Value = 'clickWithKeys(locator,keys)'
$ws.Cells.Item(56, 25).Value = 'close()'
$ws.Cells.Item(57, 25).Value = 'closeAll()'
$ws.Cells.Item(58, 25).Value = 'deselect(locator,text)'
$ws.Cells.Item(59, 25).Value = 'deselectMulti(locator,array)'
$ws.Cells.Item(60, 25).Value = 'dismissInvalidCert()'
$ws.Cells.Item(61, 25).Value = 'dismissInvalidCertPopup()'
$ws.Cells.Item(62, 25).Value = 'doubleClick(locator)'
$ws.Cells.Item(63, 25).Value = 'doubleClickAndWait(locator,waitMs)'
$ws.Cells.Item(64, 25).Value = 'doubleClickByLabel(label)'
$ws.Cells.Item(65, 25).Value = 'doubleClickByLabelAndWait(label,waitMs)'
$ws.Cells.Item(66, 25).Value = 'dragAndDrop(fromLocator,toLocator)'
$ws.Cells.Item(67, 25).Value = 'dragTo(fromLocator,xOffset,yOffset)'
$ws.Cells.Item(68, 25).Value = 'editLocalStorage(key,value)'
$ws.Cells.Item(69, 25).Value = 'executeScript(var,script)'
$ws.Cells.Item(70, 25).Value = 'focus(locator)'
$ws.Cells.Item(71, 25).Value = 'goBack()'
$ws.Cells.Item(72, 25).Value = 'goBackAndWait()'
$ws.Cells.Item(73, 25).Value = 'maximizeWindow()'
$ws.Cells.Item(74, 25).Value = 'mouseOver(locator)'
$ws.Cells.Item(75, 25).Value = 'open(url)'
$ws.Cells.Item(76, 25).Value = 'openAndWait(url,waitMs)'
$ws.Cells.Item(77, 25).Value = 'openHttpBasic(url,username,password)'
$ws.Cells.Item(78, 25).Value = 'openIgnoreTimeout(url)'
$ws.Cells.Item(79, 25).Value = 'refresh()'
$ws.Cells.Item(80, 25).Value = 'refreshAndWait()'
$ws.Cells.Item(81, 25).Value = 'resizeWindow(width,height)'
$ws.Cells.Item(82, 25).Value = 'rightClick(locator)'
$ws.Cells.Item(83, 25).Value = 'saveAllWindowIds(var)'
$ws.Cells.Item(84, 25).Value = 'saveAllWindowNames(var)'
$ws.Cells.Item(85, 25).Value = 'saveAttribute(var,locator,attrName)'
$ws.Cells.Item(86, 25).Value = 'saveAttributeList(var,locator,attrName)'
$ws.Cells.Item(87, 25).Value = 'saveCount(var,locator)'
$ws.Cells.Item(88, 25).Value = 'saveDivsAsCsv(headers,rows,cells,nextPage,file)'
$ws.Cells.Item(89, 25).Value = 'saveElement(var,locator)'
$ws.Cells.Item(90, 25).Value = 'saveElements(var,locator)'
$ws.Cells.Item(91, 25).Value = 'saveLocalStorage(var,key)'
$ws.Cells.Item(92, 25).Value = 'saveLocation(var)'
$ws.Cells.Item(93, 25).Value = 'savePageAs(var,sessionIdName,url)'
$ws.Cells.Item(94, 25).Value = 'savePageAsFile(sessionIdName,url,file)'
$ws.Cells.Item(95, 25).Value = 'saveTableAsCsv(locator,nextPageLocator,file)'
$ws.Cells.Item(96, 25).Value = 'saveText(var,locator)'
$ws.Cells.Item(97, 25).Value = 'saveTextArray(var,locator)'
$ws.Cells.Item(98, 25).Value = 'saveTextSubstringAfter(var,locator,delim)'
$ws.Cells.Item(99, 25).Value = 'saveTextSubstringBefore(var,locator,delim)'
$ws.Cells.Item(100, 25).Value = 'saveTextSubstringBetween(var,locator,start,end)'
$ws.Cells.Item(101, 25).Value = 'saveValue(var,locator)'
$ws.Cells.Item(102, 25).Value = 'saveValues(var,locator)'
$ws.Cells.Item(103, 25).Value = 'scrollElement(locator,xOffset,yOffset)'
$ws.Cells.Item(104, 25).Value = 'scrollLeft(locator,pixel)'
$ws.Cells.Item(105, 25).Value = 'scrollPage(xOffset,yOffset)'
$ws.Cells.Item(106, 25).Value = 'scrollRight(locator,pixel)'
$ws.Cells.Item(107, 25).Value = 'scrollTo(locator)'
$ws.Cells.Item(108, 25).Value = 'select(locator,text)'
$ws.Cells.Item(109, 25).Value = 'selectFrame(locator)'
$ws.Cells.Item(110, 25).Value = 'selectMulti(locator,array)'
$ws.Cells.Item(111, 25).Value = 'selectMultiOptions(locator)'
$ws.Cells.Item(112, 25).Value = 'selectText(locator)'
$ws.Cells.Item(113, 25).Value = 'selectWindow(winId)'
$ws.Cells.Item(114, 25).Value = 'selectWindowAndWait(winId,waitMs)'
$ws.Cells.Item(115, 25).Value = 'selectWindowByIndex(index)'
$ws.Cells.Item(116, 25).Value = 'selectWindowByIndexAndWait(index,waitMs)'
$ws.Cells.Item(117, 25).Value = 'toggleSelections(locator)'
$ws.Cells.Item(118, 25).Value = 'type(locator,value)'
$ws.Cells.Item(119, 25).Value = 'typeKeys(locator,value)'
$ws.Cells.Item(120, 25).Value = 'uncheckAll(locator)'
$ws.Cells.Item(121, 25).Value = 'unselectAllText()'
$ws.Cells.Item(122, 25).Value = 'upload(fieldLocator,file)'
$ws.Cells.Item(123, 25).Value = 'verifyContainText(locator,text)'
$ws.Cells.Item(124, 25).Value = 'verifyText(locator,text)'
$ws.Cells.Item(125, 25).Value = 'wait(waitMs)'
$ws.Cells.Item(126, 25).Value = 'waitForElementPresent(locator)'
$ws.Cells.Item(127, 25).Value = 'waitForPopUp(winId,waitMs)'
$ws.Cells.Item(128, 25).Value = 'waitForTextPresent(text)'
$ws.Cells.Item(129, 25).Value = 'waitForTitle(text)'

# ---- Clear now-unused trailing cells ----
$ws.Cells.Item(1, 31).ClearContents()
$ws.Cells.Item(2, 31).ClearContents()
$ws.Cells.Item(3, 31).ClearContents()
$ws.Cells.Item(4, 31).ClearContents()
$ws.Cells.Item(5, 31).ClearContents()
$ws.Cells.Item(6, 31).ClearContents()
$ws.Cells.Item(7, 31).ClearContents()
$ws.Cells.Item(8, 31).ClearContents()
$ws.Cells.Item(9, 29).ClearContents()
$ws.Cells.Item(9, 31).ClearContents()
$ws.Cells.Item(9, 26).ClearContents()
$ws.Cells.Item(10, 29).ClearContents()
$ws.Cells.Item(10, 31).ClearContents()
$ws.Cells.Item(10, 26).ClearContents()
$ws.Cells.Item(11, 29).ClearContents()
$ws.Cells.Item(11, 31).ClearContents()
$ws.Cells.Item(11, 26).ClearContents()
$ws.Cells.Item(12, 29).ClearContents()
$ws.Cells.Item(12, 31).ClearContents()
$ws.Cells.Item(12, 26).ClearContents()
$ws.Cells.Item(13, 29).ClearContents()
$ws.Cells.Item(13, 31).ClearContents()
$ws.Cells.Item(13, 26).ClearContents()
$ws.Cells.Item(14, 29).ClearContents()
$ws.Cells.Item(14, 31).ClearContents()
$ws.Cells.Item(14, 26).ClearContents()
$ws.Cells.Item(15, 29).ClearContents()
$ws.Cells.Item(15, 31).ClearContents()
$ws.Cells.Item(15, 26).ClearContents()
$ws.Cells.Item(16, 29).ClearContents()
$ws.Cells.Item(16, 31).ClearContents()
$ws.Cells.Item(16, 26).ClearContents()
$ws.Cells.Item(17, 29).ClearContents()
$ws.Cells.Item(17, 31).ClearContents()
$ws.Cells.Item(17, 26).ClearContents()
$ws.Cells.Item(18, 31).ClearContents()
$ws.Cells.Item(18, 26).ClearContents()
$ws.Cells.Item(19, 31).ClearContents()
$ws.Cells.Item(19, 26).ClearContents()
$ws.Cells.Item(20, 31).ClearContents()
$ws.Cells.Item(20, 26).ClearContents()
$ws.Cells.Item(21, 31).ClearContents()
$ws.Cells.Item(21, 26).ClearContents()
$ws.Cells.Item(22, 31).ClearContents()
$ws.Cells.Item(22, 26).ClearContents()
$ws.Cells.Item(23, 31).ClearContents()
$ws.Cells.Item(23, 26).ClearContents()
$ws.Cells.Item(24, 31).ClearContents()
$ws.Cells.Item(24, 26).ClearContents()
$ws.Cells.Item(25, 31).ClearContents()
$ws.Cells.Item(25, 26).ClearContents()
$ws.Cells.Item(26, 31).ClearContents()
$ws.Cells.Item(26, 26).ClearContents()
$ws.Cells.Item(27, 31).ClearContents()
$ws.Cells.Item(27, 26).ClearContents()
$ws.Cells.Item(28, 26).ClearContents()
$ws.Cells.Item(29, 26).ClearContents()
$ws.Cells.Item(30, 26).ClearContents()
$ws.Cells.Item(31, 1).ClearContents()
$ws.Cells.Item(31, 26).ClearContents()
$ws.Cells.Item(32, 26).ClearContents()
$ws.Cells.Item(33, 26).ClearContents()
$ws.Cells.Item(34, 26).ClearContents()
$ws.Cells.Item(35, 26).ClearContents()
$ws.Cells.Item(36, 26).ClearContents()
$ws.Cells.Item(37, 26).ClearContents()
$ws.Cells.Item(38, 26).ClearContents()
$ws.Cells.Item(39, 26).ClearContents()
$ws.Cells.Item(40, 26).ClearContents()
$ws.Cells.Item(41, 26).ClearContents()
$ws.Cells.Item(42, 26).ClearContents()
$ws.Cells.Item(43, 26).ClearContents()
$ws.Cells.Item(44, 26).ClearContents()
$ws.Cells.Item(45, 26).ClearContents()
$ws.Cells.Item(46, 26).ClearContents()
$ws.Cells.Item(47, 26).ClearContents()
$ws.Cells.Item(48, 26).ClearContents()
$ws.Cells.Item(49, 26).ClearContents()
$ws.Cells.Item(50, 26).ClearContents()
$ws.Cells.Item(51, 26).ClearContents()
$ws.Cells.Item(52, 26).ClearContents()
$ws.Cells.Item(53, 26).ClearContents()
$ws.Cells.Item(54, 26).ClearContents()
$ws.Cells.Item(55, 26).ClearContents()
$ws.Cells.Item(56, 26).ClearContents()
$ws.Cells.Item(57, 26).ClearContents()
$ws.Cells.Item(58, 26).ClearContents()
$ws.Cells.Item(59, 26).ClearContents()
$ws.Cells.Item(60, 26).ClearContents()
$ws.Cells.Item(61, 26).ClearContents()
$ws.Cells.Item(62, 26).ClearContents()
$ws.Cells.Item(63, 26).ClearContents()
$ws.Cells.Item(64, 26).ClearContents()
$ws.Cells.Item(65, 26).ClearContents()
$ws.Cells.Item(66, 26).ClearContents()
$ws.Cells.Item(67, 26).ClearContents()
$ws.Cells.Item(68, 26).ClearContents()
$ws.Cells.Item(69, 26).ClearContents()
$ws.Cells.Item(70, 26).ClearContents()
$ws.Cells.Item(71, 26).ClearContents()
$ws.Cells.Item(72, 26).ClearContents()
$ws.Cells.Item(73, 26).ClearContents()
$ws.Cells.Item(74, 26).ClearContents()
$ws.Cells.Item(75, 26).ClearContents()
$ws.Cells.Item(76, 26).ClearContents()
$ws.Cells.Item(77, 26).ClearContents()
$ws.Cells.Item(78, 26).ClearContents()
$ws.Cells.Item(79, 26).ClearContents()
$ws.Cells.Item(80, 26).ClearContents()
$ws.Cells.Item(81, 26).ClearContents()
$ws.Cells.Item(82, 26).ClearContents()
$ws.Cells.Item(83, 26).ClearContents()
$ws.Cells.Item(84, 26).ClearContents()
$ws.Cells.Item(85, 26).ClearContents()
$ws.Cells.Item(86, 26).ClearContents()
$ws.Cells.Item(87, 26).ClearContents()
$ws.Cells.Item(88, 26).ClearContents()
$ws.Cells.Item(89, 26).ClearContents()
$ws.Cells.Item(90, 26).ClearContents()
$ws.Cells.Item(91, 26).ClearContents()
$ws.Cells.Item(92, 26).ClearContents()
$ws.Cells.Item(93, 26).ClearContents()
$ws.Cells.Item(94, 26).ClearContents()
$ws.Cells.Item(95, 26).ClearContents()
$ws.Cells.Item(96, 26).ClearContents()
$ws.Cells.Item(97, 26).ClearContents()
$ws.Cells.Item(98, 26).ClearContents()
$ws.Cells.Item(99, 26).ClearContents()
$ws.Cells.Item(100, 26).ClearContents()
$ws.Cells.Item(101, 26).ClearContents()
$ws.Cells.Item(102, 26).ClearContents()
$ws.Cells.Item(103, 26).ClearContents()
$ws.Cells.Item(104, 26).ClearContents()
$ws.Cells.Item(105, 26).ClearContents()
$ws.Cells.Item(106, 26).ClearContents()
$ws.Cells.Item(107, 26).ClearContents()
$ws.Cells.Item(108, 26).ClearContents()
$ws.Cells.Item(109, 26).ClearContents()
$ws.Cells.Item(110, 26).ClearContents()
$ws.Cells.Item(111, 26).ClearContents()
$ws.Cells.Item(112, 26).ClearContents()
$ws.Cells.Item(113, 26).ClearContents()
$ws.Cells.Item(114, 26).ClearContents()
$ws.Cells.Item(115, 26).ClearContents()
$ws.Cells.Item(116, 26).ClearContents()
$ws.Cells.Item(117, 26).ClearContents()
$ws.Cells.Item(118, 26).ClearContents()
$ws.Cells.Item(119, 26).ClearContents()
$ws.Cells.Item(120, 26).ClearContents()
$ws.Cells.Item(121, 26).ClearContents()
$ws.Cells.Item(122, 26).ClearContents()
$ws.Cells.Item(123, 26).ClearContents()
$ws.Cells.Item(124, 26).ClearContents()
$ws.Cells.Item(125, 26).ClearContents()
$ws.Cells.Item(126, 26).ClearContents()
$ws.Cells.Item(127, 26).ClearContents()
$ws.Cells.Item(128, 26).ClearContents()
$ws.Cells.Item(129, 26).ClearContents()

# ---- Update defined name ranges affected by the column/row shift ----
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
